# Auto-applies the cryptos.xlsx data refresh described in the commit
# 'Updated cryptos list on Thu Aug 10 12:12:59 UTC 2023 with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.475.48"
$ws.Range("E2").Value = "  -1.38%  "

# Row 3
$ws.Range("D3").Value = "1.849.00"
$ws.Range("E3").Value = "  -0.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'241.80"
$ws.Range("E5").Value = "  -0.95%  "

# Row 6
$ws.Range("D6").Value = "'0.6291"
$ws.Range("E6").Value = "  -2.67%  "

# Row 8
$ws.Range("D8").Value = "'0.07515"
$ws.Range("E8").Value = "  -0.46%  "

# Row 9
$ws.Range("D9").Value = "'0.2972"
$ws.Range("E9").Value = "  -0.56%  "

# Row 10
$ws.Range("D10").Value = "'24.36"
$ws.Range("E10").Value = "  -2.13%  "

# Row 11
$ws.Range("D11").Value = "'0.07735"

# Row 12
$ws.Range("D12").Value = "1.884.63"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("D13").Value = "'0.6923"
$ws.Range("E13").Value = "  -0.19%  "

# Row 14
$ws.Range("D14").Value = "'5.002"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("D15").Value = "'83.43"
$ws.Range("E15").Value = "  -0.61%  "

# Row 16
$ws.Range("D16").Value = "'0.000009799"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17
$ws.Range("D17").Value = "2.114.89"
$ws.Range("E17").Value = "  -0.68%  "

# Row 18
$ws.Range("D18").Value = "'6.236"
$ws.Range("E18").Value = "  +1.74%  "

# Row 19
$ws.Range("D19").Value = "29.507.30"
$ws.Range("E19").Value = "  -1.30%  "

# Row 20
$ws.Range("D20").Value = "'232.66"
$ws.Range("E20").Value = "  -1.70%  "

# Row 21
$ws.Range("D21").Value = "'12.50"
$ws.Range("E21").Value = "  -1.49%  "

# Row 22
$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("D23").Value = "'7.630"
$ws.Range("E23").Value = "  -0.80%  "

# Row 25
$ws.Range("D25").Value = "'154.93"
$ws.Range("E25").Value = "  -2.50%  "

# Row 26
$ws.Range("D26").Value = "'0.1388"
$ws.Range("E26").Value = "  -2.64%  "

# Row 27
$ws.Range("D27").Value = "'8.447"
$ws.Range("E27").Value = "  -1.50%  "

# Row 28
$ws.Range("D28").Value = "'17.69"
$ws.Range("E28").Value = "  -1.47%  "

# Row 29
$ws.Range("D29").Value = "'1.474"
$ws.Range("E29").Value = "  -1.49%  "

# Row 30
$ws.Range("D30").Value = "'0.05932"
$ws.Range("E30").Value = "  -4.56%  "

# Row 31
$ws.Range("D31").Value = "'1.253"
$ws.Range("E31").Value = "  -2.69%  "

# Row 32
$ws.Range("D32").Value = "'4.106"
$ws.Range("E32").Value = "  -1.48%  "

# Row 33
$ws.Range("D33").Value = "'4.021"
$ws.Range("E33").Value = "  -1.92%  "

# Row 34
$ws.Range("D34").Value = "'1.872"
$ws.Range("E34").Value = "  -0.82%  "

# Row 35
$ws.Range("D35").Value = "'1.165"
$ws.Range("E35").Value = "  -0.51%  "

# Row 36
$ws.Range("D36").Value = "'0.7178"
$ws.Range("E36").Value = "  -2.06%  "

# Row 37
$ws.Range("D37").Value = "'2.589"
$ws.Range("E37").Value = "  -1.27%  "

# Row 38
$ws.Range("D38").Value = "'2.799"
$ws.Range("E38").Value = "  -0.61%  "

# Row 39
$ws.Range("D39").Value = "1.237.06"
$ws.Range("E39").Value = "  +1.64%  "

# Row 40
$ws.Range("D40").Value = "'0.01794"
$ws.Range("E40").Value = "  +0.39%  "

# Row 41
$ws.Range("D41").Value = "'0.9068"
$ws.Range("E41").Value = "  -1.33%  "

# Row 42
$ws.Range("D42").Value = "'6.089"
$ws.Range("E42").Value = "  -4.80%  "

# Row 43
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "2.036.28"
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("E45").Value = "  -0.83%  "

# Row 46
$ws.Range("D46").Value = "'67.13"
$ws.Range("E46").Value = "  -0.48%  "

# Row 47
$ws.Range("D47").Value = "'7.340"
$ws.Range("E47").Value = "  +8.99%  "

# Row 48
$ws.Range("D48").Value = "'0.00000000118"
$ws.Range("E48").Value = "  -1.99%  "

# Row 49
$ws.Range("D49").Value = "'0.4034"
$ws.Range("E49").Value = "  -1.04%  "

# Row 50
$ws.Range("D50").Value = "'9.121"
$ws.Range("E50").Value = "  -1.00%  "

# Row 51
$ws.Range("D51").Value = "'1.702"
$ws.Range("E51").Value = "  +1.91%  "
